$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H106").Value = 4309.125
$ws.Range("I106").Value = 3579
$ws.Range("K106").Value = 3579
$ws.Range("M106").Value = -2948
$ws.Range("H112").Value = 1680.9362
$ws.Range("I112").Value = 1050
$ws.Range("J112").Value = 1708.9778
$ws.Range("K112").Value = 3150
$ws.Range("L112").Value = 5126.9334
$ws.Range("M112").Value = -2042
$ws.Range("N112").Value = -7342.9334
$ws.Range("H132").Value = 6044.088
$ws.Range("I132").Value = 6469.42
$ws.Range("K132").Value = 19408.26
$ws.Range("M132").Value = -16878.26
$ws.Range("H138").Value = 4517.271
$ws.Range("I138").Value = 1663.9166
$ws.Range("J138").Value = 5468.3887
$ws.Range("K138").Value = 4991.7498
$ws.Range("L138").Value = 16405.1661
$ws.Range("M138").Value = 148.2502000000004
$ws.Range("N138").Value = -26685.1661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3701.861
$ws.Range("I45").Value = 3959.3572
$ws.Range("J45").Value = 2800.625
$ws.Range("K45").Value = 3959.3572
$ws.Range("L45").Value = 2800.625
$ws.Range("M45").Value = -3582.3572
$ws.Range("N45").Value = -3554.625
$ws.Range("H46").Value = 12384
$ws.Range("J46").Value = 12384
$ws.Range("L46").Value = 12384
$ws.Range("N46").Value = -13022
$ws.Range("H61").Value = 3939.9167
$ws.Range("I61").Value = 3227.95
$ws.Range("K61").Value = 3227.95
$ws.Range("M61").Value = -3015.95
$ws.Range("H63").Value = 7464.4546
$ws.Range("I63").Value = 4576.5
$ws.Range("K63").Value = 4576.5
$ws.Range("M63").Value = -3890.5
$ws.Range("H66").Value = 7464.4546
$ws.Range("I66").Value = 4576.5
$ws.Range("K66").Value = 22882.5
$ws.Range("M66").Value = -19450.5
$ws.Range("H74").Value = 2637.8823
$ws.Range("I74").Value = 1231.56
$ws.Range("J74").Value = 6544.3335
$ws.Range("K74").Value = 1231.56
$ws.Range("L74").Value = 6544.3335
$ws.Range("M74").Value = -357.5599999999999
$ws.Range("N74").Value = -8292.333500000001
$ws.Range("H77").Value = 2637.8823
$ws.Range("I77").Value = 1231.56
$ws.Range("J77").Value = 6544.3335
$ws.Range("K77").Value = 6157.799999999999
$ws.Range("L77").Value = 32721.6675
$ws.Range("M77").Value = -1789.799999999999
$ws.Range("N77").Value = -41457.6675
$ws.Range("H122").Value = 2885.4194
$ws.Range("I122").Value = 2075.15
$ws.Range("J122").Value = 4358.636
$ws.Range("K122").Value = 6225.450000000001
$ws.Range("L122").Value = 13075.908
$ws.Range("M122").Value = -3775.450000000001
$ws.Range("N122").Value = -17975.908
$ws.Range("H132").Value = 3746.3333
$ws.Range("I132").Value = 3799.6428
$ws.Range("K132").Value = 11398.9284
$ws.Range("M132").Value = -8868.928400000001
$ws.Range("H136").Value = 3939.9167
$ws.Range("I136").Value = 3227.95
$ws.Range("K136").Value = 9683.849999999999
$ws.Range("M136").Value = -7133.849999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1812.3334
$ws.Range("I36").Value = 1812.3334
$ws.Range("K36").Value = 1812.3334
$ws.Range("M36").Value = -1278.3334
$ws.Range("H64").Value = 896.5
$ws.Range("I64").Value = 149
$ws.Range("J64").Value = 1145.6666
$ws.Range("K64").Value = 149
$ws.Range("L64").Value = 1145.6666
$ws.Range("M64").Value = 76
$ws.Range("N64").Value = -1595.6666
$ws.Range("H67").Value = 896.5
$ws.Range("I67").Value = 149
$ws.Range("J67").Value = 1145.6666
$ws.Range("K67").Value = 149
$ws.Range("L67").Value = 1145.6666
$ws.Range("M67").Value = 631
$ws.Range("N67").Value = -2705.6666
$ws.Range("H94").Value = 642.8095
$ws.Range("I94").Value = 477
$ws.Range("K94").Value = 477
$ws.Range("M94").Value = -26
$ws.Range("H99").Value = 39367.5
$ws.Range("I99").Value = 60427.555
$ws.Range("J99").Value = 1459.4
$ws.Range("K99").Value = 60427.555
$ws.Range("L99").Value = 1459.4
$ws.Range("M99").Value = -58929.555
$ws.Range("N99").Value = -4455.4
$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8452.081
$ws.Range("I31").Value = 875.8570999999999
$ws.Range("J31").Value = 24362.15
$ws.Range("K31").Value = 875.8570999999999
$ws.Range("L31").Value = 24362.15
$ws.Range("M31").Value = -580.8570999999999
$ws.Range("N31").Value = -24952.15
$ws.Range("H34").Value = 8452.081
$ws.Range("I34").Value = 875.8570999999999
$ws.Range("J34").Value = 24362.15
$ws.Range("K34").Value = 875.8570999999999
$ws.Range("L34").Value = 24362.15
$ws.Range("M34").Value = -673.8570999999999
$ws.Range("N34").Value = -24766.15
$ws.Range("H58").Value = 3846.9443
$ws.Range("I58").Value = 2995.0908
$ws.Range("J58").Value = 5185.5713
$ws.Range("K58").Value = 2995.0908
$ws.Range("L58").Value = 5185.5713
$ws.Range("M58").Value = -2792.0908
$ws.Range("N58").Value = -5591.5713
$ws.Range("H62").Value = 2700
$ws.Range("I62").Value = 2366.6667
$ws.Range("K62").Value = 2366.6667
$ws.Range("M62").Value = -1742.6667
$ws.Range("H65").Value = 2700
$ws.Range("I65").Value = 2366.6667
$ws.Range("K65").Value = 11833.3335
$ws.Range("M65").Value = -8713.333500000001
$ws.Range("H132").Value = 1421.1111
$ws.Range("I132").Value = 1438.0769
$ws.Range("J132").Value = 1377
$ws.Range("K132").Value = 4314.2307
$ws.Range("L132").Value = 4131
$ws.Range("M132").Value = -1784.2307
$ws.Range("N132").Value = -9191
$ws.Range("H134").Value = 26002.143
$ws.Range("I134").Value = 27871.54
$ws.Range("K134").Value = 83614.62
$ws.Range("M134").Value = -81079.62
$ws.Range("H136").Value = 3846.9443
$ws.Range("I136").Value = 2995.0908
$ws.Range("J136").Value = 5185.5713
$ws.Range("K136").Value = 8985.2724
$ws.Range("L136").Value = 15556.7139
$ws.Range("M136").Value = -6435.2724
$ws.Range("N136").Value = -20656.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 55901290
$ws.Range("I4").Value = 85100090
$ws.Range("K4").Value = 255300270
$ws.Range("M4").Value = -255300158
$ws.Range("H5").Value = 2744.5715
$ws.Range("I5").Value = 802.6667
$ws.Range("K5").Value = 2408.0001
$ws.Range("M5").Value = -2296.0001
$ws.Range("H111").Value = 4268.8
$ws.Range("I111").Value = 2586
$ws.Range("K111").Value = 7758
$ws.Range("M111").Value = -4691
$ws.Range("H122").Value = 11765125
$ws.Range("J122").Value = 22222588
$ws.Range("L122").Value = 200003292
$ws.Range("N122").Value = -200008192
$ws.Range("H135").Value = 2744.5715
$ws.Range("I135").Value = 802.6667
$ws.Range("K135").Value = 7224.0003
$ws.Range("M135").Value = -4689.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2842.1072
$ws.Range("I102").Value = 2164.4614
$ws.Range("J102").Value = 3429.4
$ws.Range("K102").Value = 2164.4614
$ws.Range("L102").Value = 3429.4
$ws.Range("M102").Value = -542.4614000000001
$ws.Range("N102").Value = -6673.4
$ws.Range("H113").Value = 4606.2583
$ws.Range("I113").Value = 4762.2383
$ws.Range("J113").Value = 4278.7
$ws.Range("K113").Value = 4762.2383
$ws.Range("L113").Value = 4278.7
$ws.Range("M113").Value = -2592.2383
$ws.Range("N113").Value = -8618.700000000001
$ws.Range("H122").Value = 4453.381
$ws.Range("I122").Value = 7252.5557
$ws.Range("J122").Value = 2354
$ws.Range("K122").Value = 21757.6671
$ws.Range("L122").Value = 7062
$ws.Range("M122").Value = -19307.6671
$ws.Range("N122").Value = -11962
$ws.Range("H126").Value = 3967.75
$ws.Range("I126").Value = 3259
$ws.Range("K126").Value = 9777
$ws.Range("M126").Value = -7307
$ws.Range("H132").Value = 3747.709
$ws.Range("I132").Value = 3076.2666
$ws.Range("K132").Value = 9228.799800000001
$ws.Range("M132").Value = -6698.799800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3262.818
$ws.Range("I7").Value = 3136.6875
$ws.Range("J7").Value = 3599.1667
$ws.Range("K7").Value = 3136.6875
$ws.Range("L7").Value = 3599.1667
$ws.Range("M7").Value = -3024.6875
$ws.Range("N7").Value = -3823.1667
$ws.Range("H40").Value = 4780.615
$ws.Range("J40").Value = 4499.5
$ws.Range("L40").Value = 4499.5
$ws.Range("N40").Value = -4771.5
$ws.Range("H126").Value = 3262.818
$ws.Range("I126").Value = 3136.6875
$ws.Range("J126").Value = 3599.1667
$ws.Range("K126").Value = 9410.0625
$ws.Range("L126").Value = 10797.5001
$ws.Range("M126").Value = -6940.0625
$ws.Range("N126").Value = -15737.5001
$ws.Range("H132").Value = 2722.2
$ws.Range("I132").Value = 2540.2415
$ws.Range("J132").Value = 7999
$ws.Range("K132").Value = 7620.7245
$ws.Range("L132").Value = 23997
$ws.Range("M132").Value = -5090.7245
$ws.Range("N132").Value = -29057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5163.1665
$ws.Range("I107").Value = 2499.5
$ws.Range("K107").Value = 7498.5
$ws.Range("M107").Value = -5578.5
$ws.Range("H126").Value = 3488.6667
$ws.Range("I126").Value = 2837.7693
$ws.Range("K126").Value = 8513.3079
$ws.Range("M126").Value = -6043.3079
$ws.Range("H132").Value = 5242.524
$ws.Range("I132").Value = 1399.3846
$ws.Range("J132").Value = 11487.625
$ws.Range("K132").Value = 4198.1538
$ws.Range("L132").Value = 34462.875
$ws.Range("M132").Value = -1668.1538
$ws.Range("N132").Value = -39522.875
$ws.Range("H136").Value = 30380.363
$ws.Range("I136").Value = 82079.336
$ws.Range("K136").Value = 246238.008
$ws.Range("M136").Value = -243688.008
$ws.Range("H138").Value = 77857.25
$ws.Range("J138").Value = 77857.25
$ws.Range("L138").Value = 77857.25
$ws.Range("N138").Value = -88137.25
